# Apply the "cryptos list" refresh described in the commit:
#   "Updated cryptos list on Fri Oct 18 08:17:59 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''67.885.44'
$ws.Range('E2').Value = '  +1.00%  '

$ws.Range('D3').Value = '''2.638.26'
$ws.Range('E3').Value = '  +0.63%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '''598.07'
$ws.Range('E5').Value = '  +0.33%  '

$ws.Range('D6').Value = '''153.71'
$ws.Range('E6').Value = '  +0.98%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('E8').Value = '  -0.50%  '

$ws.Range('D9').Value = '''2.636.09'
$ws.Range('E9').Value = '  +0.58%  '

$ws.Range('E10').Value = '  +10.84%  '

$ws.Range('E11').Value = '  -0.64%  '

$ws.Range('E12').Value = '  +0.74%  '

$ws.Range('E13').Value = '  +0.00%  '

$ws.Range('D14').Value = '''27.64'
$ws.Range('E14').Value = '  +0.44%  '

$ws.Range('E15').Value = '  +3.92%  '

$ws.Range('D16').Value = '''3.119.67'
$ws.Range('E16').Value = '  +0.51%  '

$ws.Range('D17').Value = '''67.738.61'

$ws.Range('D18').Value = '''2.650.19'
$ws.Range('E18').Value = '  +1.19%  '

$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '''11.43'
$ws.Range('E19').Value = '  +3.00%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '''372.76'
$ws.Range('E20').Value = '  +2.82%  '

$ws.Range('D21').Value = '''7.47'
$ws.Range('E21').Value = '  -0.09%  '

$ws.Range('D22').Value = '''4.26'
$ws.Range('E22').Value = '  -0.80%  '

$ws.Range('E23').Value = '  -1.29%  '

$ws.Range('D24').Value = '''2.05'
$ws.Range('E24').Value = '  -2.25%  '

$ws.Range('D25').Value = '''72.15'
$ws.Range('E25').Value = '  +1.76%  '

$ws.Range('E26').Value = '  +0.10%  '

$ws.Range('D27').Value = '''9.96'
$ws.Range('E27').Value = '  -0.92%  '

$ws.Range('D28').Value = '''2.760.30'

$ws.Range('E29').Value = '  +2.46%  '

$ws.Range('E30').Value = '  +0.06%  '

$ws.Range('D31').Value = '''574.09'
$ws.Range('E31').Value = '  -0.67%  '

$ws.Range('E32').Value = '  +1.05%  '

$ws.Range('D33').Value = '''7.88'
$ws.Range('E33').Value = '  +1.28%  '

$ws.Range('D34').Value = '''1.84'
$ws.Range('E34').Value = '  +0.59%  '

$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.05%  '

$ws.Range('D36').Value = '''0.125'
$ws.Range('E36').Value = '  -1.57%  '

$ws.Range('E37').Value = '  +0.33%  '

$ws.Range('D38').Value = '''158.84'
$ws.Range('E38').Value = '  +1.20%  '

$ws.Range('D39').Value = '''19.17'
$ws.Range('E39').Value = '  +0.39%  '

$ws.Range('E40').Value = '  +5.69%  '

$ws.Range('E41').Value = '  +0.52%  '

$ws.Range('D42').Value = '''5.36'
$ws.Range('E42').Value = '  +2.36%  '

$ws.Range('E43').Value = '  +2.82%  '

$ws.Range('D44').Value = '''0.0₆0324'
$ws.Range('E44').Value = '  +14.90%  '

$ws.Range('E45').Value = '  +4.79%  '

$ws.Range('E46').Value = '  +0.01%  '

$ws.Range('D47').Value = '''40.33'
$ws.Range('E47').Value = '  -1.96%  '

$ws.Range('D48').Value = '''155.74'
$ws.Range('E48').Value = '  -0.28%  '

$ws.Range('E49').Value = '  -0.53%  '

$ws.Range('D50').Value = '''22.30'
$ws.Range('E50').Value = '  +8.78%  '

$ws.Range('E51').Value = '  -1.29%  '
